$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = -0.6141830484198026
$ws.Range("J5").Value = 0.451704197853215
$ws.Range("K5").Value = 0.1321754508754534
$ws.Range("L5").Value = 2.526671512074564
